$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Locate the existing blank paragraph that currently sits right after
# "Author: Adam Foster" (and right before "US - replicated paper
# findings..."). We will keep this paragraph completely untouched and
# instead insert a second, brand-new blank paragraph immediately
# before it, so that two empty paragraphs remain between our new
# content and the "US - replicated..." paragraph.
# ------------------------------------------------------------------
$authorRange = $d.Content.Find.Execute("Author: Adam Foster", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

$authorPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.StartsWith("Author: Adam Foster")) {
        $authorPara = $p
        break
    }
}

$existingBlank = $authorPara.Next()
$blankStart = $existingBlank.Range
$blankStart.Collapse(1)
$blankStart.InsertParagraphBefore()

# ------------------------------------------------------------------
# Insert the new body paragraphs (leading blank line + three text
# paragraphs) right after "Author: Adam Foster".
# ------------------------------------------------------------------
$r = $authorPara.Range
$r.Collapse(0)
$r.InsertAfter("`rThe goal of this project was to test the market timing strategy described in Asness, C., Ilmanen, A., & Maloney, T. (2017). Market timing: Sin a little resolving the valuation timing puzzle. Journal of Investment Management, 15(3), 23-40 on a chosen.`rThe German equity market was the market of choice with sufficient data available. The analysis was produced using the cyclically-adjusted P/E ratio (CAPE) for Germany created by Barclays Research and DAX index prices on Yahoo Finance. The reason for the divergence in data sources was the lack of index price data in the Barclays Research dataset, unlike the completeness of US data produced by Robert Shiller. Yahoo Finance contained one of the longest uninterrupted DAX time series which was deemed an acceptable alternative and likely to overlap with the German equity selection used in CAPE. CAPE data extended from Jan-82 to May-23 and DAX data was from Jan-88 to Apr-23, thus the time frame considered spanned across the common Jan-88 to Apr-23 period. This is several decades shorter than the US data provided by Shiller, but still worthy of analysis – long-term returns were reduced from 10 years in the paper to 5 years and the rolling window of CAPE quintiles from up to 60 years to up to 20 years, still representing roughly half of the overall time period.`rInitially, US results in the paper were reproduced.")

# ------------------------------------------------------------------
# Italicise the journal-article citation inside the first new
# paragraph (matches the <w:i/><w:iCs/> run in the target document).
# ------------------------------------------------------------------
$citationRange = $d.Content
$citationRange.Find.ClearFormatting()
$citationRange.Find.Execute("Asness, C., Ilmanen, A., & Maloney, T. (2017). Market timing: Sin a little resolving the valuation timing puzzle. Journal of Investment Management, 15(3), 23-40", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($citationRange.Find.Found) {
    $citationRange.Font.Italic = 1
}
